$d = $word.ActiveDocument

$replacements = @(
    @{old="51×68="; new="50×21="},
    @{old="24×75="; new="24×70="},
    @{old="92×42="; new="93×31="},
    @{old="61×27="; new="89×80="},
    @{old="17×82="; new="68×67="},
    @{old="76×45="; new="58×58="},
    @{old="66×71="; new="45×80="},
    @{old="46×58="; new="23×21="},
    @{old="64×35="; new="55×36="},
    @{old="96×83="; new="59×40="},
    @{old="58×38="; new="71×15="},
    @{old="91×15="; new="66×39="},
    @{old="61×37="; new="98×64="},
    @{old="70×49="; new="31×81="},
    @{old="58×32="; new="18×26="},
    @{old="85×36="; new="97×68="},
    @{old="69×13="; new="95×63="},
    @{old="94×54="; new="63×32="},
    @{old="30×26="; new="73×17="},
    @{old="62×73="; new="86×28="},
    @{old="65×79="; new="45×39="},
    @{old="50×68="; new="97×79="},
    @{old="13×43="; new="31×51="},
    @{old="19×26="; new="90×64="},
    @{old="34×63="; new="34×73="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
